$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 35

# Column A holds a date-looking string ("12/29/2025"). Excel's COM layer
# auto-converts strings that look like dates into date serial numbers, so
# force the cell to Text format first, assign the literal string, then
# restore the cell's style to Normal so no stray number-format sticks
# around on the new row (matching the plain, unstyled data rows above it).
$dateCell = $ws.Cells.Item($row, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "12/29/2025"
$dateCell.Style = "Normal"

$ws.Cells.Item($row, 2).Value = 12063.57
$ws.Cells.Item($row, 3).Value = 0.212164207900799
$ws.Cells.Item($row, 4).Value = 0.787835792099201
$ws.Cells.Item($row, 5).Value = -143.07
$ws.Cells.Item($row, 6).Value = -26.74
$ws.Cells.Item($row, 7).Value = -21087.65
$ws.Cells.Item($row, 8).Value = -68.93000000000001
$ws.Cells.Item($row, 9).Value = -493.39
$ws.Cells.Item($row, 10).Value = -16.16
